# tabloların yarısı mongoDb ye aktarıldı
# Append new log rows to the "logs" sheet and a new shared error string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Id = 639; Desc = "Error"; Msg = "Tüm öğrencilerin listelenmesinde bir hata oluştu"; Date = 45618 },
    @{ Id = 640; Desc = "Error"; Msg = "Tüm öğrencilerin listelenmesinde bir hata oluştu"; Date = 45618 },
    @{ Id = 641; Desc = "Error"; Msg = "Tüm öğrencilerin listelenmesinde bir hata oluştu"; Date = 45618 },
    @{ Id = 642; Desc = "Error"; Msg = "Tüm öğrencilerin listelenmesinde bir hata oluştu"; Date = 45618 },
    @{ Id = 643; Desc = "Info";  Msg = "Tüm öğrenciler listelendi"; Date = 45618 },
    @{ Id = 644; Desc = "Succes"; Msg = "Öğrenci silme İşlemi başarılı."; Date = 45618 },
    @{ Id = 645; Desc = "Info";  Msg = "Tüm öğrenciler listelendi"; Date = 45621 }
)

$startRow = 334
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.Id
    $ws.Cells.Item($r, 2).Value = $data.Desc
    $ws.Cells.Item($r, 3).Value = $data.Msg

    $dateCell = $ws.Cells.Item($r, 4)
    $dateCell.Value = $data.Date
    $dateCell.NumberFormat = "dd-MM-yyyy"
}
